$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 98
$ws.Range("E2").Value = "Running Robot tests (Assigning) studying Assigning generously Assign all operators Assign all operators for the remaining users"
$ws.Range("F2").Value = "['Running', 'Robot', 'tests', '(', 'Assigning', ')', 'studying', 'Assigning', 'generously', 'Assign', 'all', 'operators', 'Assign', 'all', 'operators', 'for', 'the', 'remaining', 'users']"
$ws.Range("G2").Value = "['Running', 'Robot', 'tests', 'Assigning', 'studying', 'Assigning', 'generously', 'Assign', 'all', 'operators', 'Assign', 'all', 'operators', 'for', 'the', 'remaining', 'users']"
$ws.Range("H2").Value = "['Running', 'Robot', 'tests', 'Assigning', 'studying', 'Assigning', 'generously', 'Assign', 'operators', 'Assign', 'operators', 'remaining', 'users']"
$ws.Range("C3").Value = 321
$ws.Range("C4").Value = 62
$ws.Range("C5").Value = 160
$ws.Range("C6").Value = 90
$ws.Range("C7").Value = 283
$ws.Range("C8").Value = 260
$ws.Range("C9").Value = 271
$ws.Range("C10").Value = 224
$ws.Range("C11").Value = 137
$ws.Range("C12").Value = 78
$ws.Range("C13").Value = 279
$ws.Range("C14").Value = 44
$ws.Range("C15").Value = 229
$ws.Range("E15").Value = "Reports tests (Cleaning Jobs) I Me Verify `"Download Report`" button is disabled by default and the empty PDF message in the container is shown Input worksite, robot and month, then select first report shown, assert API call, download the report then verify the report has been downloaded (Excluding CC) I Me Input company, worksite, robot and month but do not select any report, verify `"Download Report`" button is disabled and correct empty PDF frame (Excluding CC) Input company, worksite, robot and month, download the first report, verify successful api call and file download user 32131 use213123r 5434634"
$ws.Range("F15").Value = "['Reports', 'tests', '(', 'Cleaning', 'Jobs', ')', 'I', 'Me', 'Verify', '````', 'Download', 'Report', `"''`", 'button', 'is', 'disabled', 'by', 'default', 'and', 'the', 'empty', 'PDF', 'message', 'in', 'the', 'container', 'is', 'shown', 'Input', 'worksite', ',', 'robot', 'and', 'month', ',', 'then', 'select', 'first', 'report', 'shown', ',', 'assert', 'API', 'call', ',', 'download', 'the', 'report', 'then', 'verify', 'the', 'report', 'has', 'been', 'downloaded', '(', 'Excluding', 'CC', ')', 'I', 'Me', 'Input', 'company', ',', 'worksite', ',', 'robot', 'and', 'month', 'but', 'do', 'not', 'select', 'any', 'report', ',', 'verify', '````', 'Download', 'Report', `"''`", 'button', 'is', 'disabled', 'and', 'correct', 'empty', 'PDF', 'frame', '(', 'Excluding', 'CC', ')', 'Input', 'company', ',', 'worksite', ',', 'robot', 'and', 'month', ',', 'download', 'the', 'first', 'report', ',', 'verify', 'successful', 'api', 'call', 'and', 'file', 'download', 'user', '32131', 'use213123r', '5434634']"
$ws.Range("G15").Value = "['Reports', 'tests', 'Cleaning', 'Jobs', 'I', 'Me', 'Verify', 'Download', 'Report', 'button', 'is', 'disabled', 'by', 'default', 'and', 'the', 'empty', 'PDF', 'message', 'in', 'the', 'container', 'is', 'shown', 'Input', 'worksite', 'robot', 'and', 'month', 'then', 'select', 'first', 'report', 'shown', 'assert', 'API', 'call', 'download', 'the', 'report', 'then', 'verify', 'the', 'report', 'has', 'been', 'downloaded', 'Excluding', 'CC', 'I', 'Me', 'Input', 'company', 'worksite', 'robot', 'and', 'month', 'but', 'do', 'not', 'select', 'any', 'report', 'verify', 'Download', 'Report', 'button', 'is', 'disabled', 'and', 'correct', 'empty', 'PDF', 'frame', 'Excluding', 'CC', 'Input', 'company', 'worksite', 'robot', 'and', 'month', 'download', 'the', 'first', 'report', 'verify', 'successful', 'api', 'call', 'and', 'file', 'download', 'user', '32131', 'use213123r', '5434634']"
$ws.Range("C16").Value = 251
$ws.Range("C17").Value = 341
$ws.Range("C18").Value = 31
$ws.Range("C19").Value = 258
$ws.Range("C20").Value = 166
$ws.Range("C21").Value = 147
$ws.Range("C22").Value = 332
$ws.Range("C23").Value = 82
$ws.Range("C24").Value = 192
$ws.Range("C25").Value = 45
$ws.Range("C26").Value = 41
$ws.Range("C27").Value = 43
$ws.Range("C28").Value = 307
$ws.Range("C29").Value = 34
$ws.Range("C30").Value = 225
$ws.Range("C31").Value = 140
$ws.Range("C32").Value = 246
$ws.Range("C33").Value = 44
$ws.Range("C34").Value = 300
$ws.Range("C35").Value = 143
$ws.Range("C36").Value = 254
$ws.Range("C37").Value = 283
$ws.Range("C38").Value = 313
$ws.Range("C39").Value = 149
$ws.Range("C40").Value = 206
$ws.Range("C41").Value = 148
$ws.Range("C42").Value = 142
$ws.Range("C43").Value = 265
$ws.Range("C44").Value = 178
$ws.Range("C45").Value = 41
$ws.Range("C46").Value = 243
$ws.Range("C47").Value = 314
$ws.Range("C48").Value = 358
$ws.Range("C49").Value = 81
$ws.Range("C50").Value = 125
$ws.Range("C51").Value = 352
$ws.Range("C52").Value = 181
$ws.Range("C53").Value = 91
$ws.Range("C54").Value = 200
$ws.Range("C55").Value = 286
$ws.Range("C56").Value = 246
$ws.Range("C57").Value = 289
$ws.Range("C58").Value = 127
$ws.Range("C59").Value = 185
$ws.Range("C60").Value = 175
$ws.Range("C61").Value = 330
$ws.Range("C62").Value = 285
$ws.Range("C63").Value = 288
$ws.Range("C64").Value = 231
$ws.Range("C65").Value = 331
$ws.Range("C66").Value = 47
$ws.Range("C67").Value = 275
$ws.Range("C68").Value = 154
$ws.Range("C69").Value = 236
$ws.Range("C70").Value = 242
$ws.Range("C71").Value = 118
$ws.Range("C72").Value = 217
$ws.Range("C73").Value = 310
$ws.Range("C74").Value = 221
$ws.Range("C75").Value = 74
$ws.Range("C76").Value = 254
$ws.Range("C77").Value = 290
$ws.Range("C78").Value = 85
$ws.Range("C79").Value = 113
$ws.Range("C80").Value = 296
$ws.Range("C81").Value = 231
$ws.Range("C82").Value = 219
$ws.Range("C83").Value = 280
